$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.363.39"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").Value = "1.874.01"
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'0.7121"
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("D6").Value = "'241.91"
$ws.Range("E6").Value = "  +0.68%  "

$ws.Range("D7").Value = "'1.000"

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.3116"
$ws.Range("E8").Value = "  +1.26%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.07795"
$ws.Range("E9").Value = "  +1.75%  "

$ws.Range("D10").Value = "'25.20"
$ws.Range("E10").Value = "  +1.88%  "

$ws.Range("D11").Value = "'0.08436"
$ws.Range("E11").Value = "  +0.77%  "

$ws.Range("D12").Value = "1.869.95"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").Value = "'5.238"
$ws.Range("E13").Value = "  +0.85%  "

$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("D15").Value = "'91.07"
$ws.Range("E15").Value = "  +0.01%  "

$ws.Range("D16").Value = "29.373.94"
$ws.Range("E16").Value = "  +0.56%  "

$ws.Range("D17").Value = "'6.070"
$ws.Range("E17").Value = "  +2.30%  "

$ws.Range("D18").Value = "'0.000008229"
$ws.Range("E18").Value = "  +5.13%  "

$ws.Range("D19").Value = "'241.01"
$ws.Range("E19").Value = "  -0.36%  "

$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("E21").Value = "  +0.30%  "

$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").Value = "'7.785"
$ws.Range("E23").Value = "  -0.62%  "

$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").Value = "'0.1595"
$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("D26").Value = "'163.19"
$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("D27").Value = "'9.072"
$ws.Range("E27").Value = "  +2.13%  "

$ws.Range("D28").Value = "'18.57"
$ws.Range("E28").Value = "  +0.63%  "

$ws.Range("D29").Value = "'1.507"
$ws.Range("E29").Value = "  +0.76%  "

$ws.Range("D30").Value = "'4.420"
$ws.Range("E30").Value = "  +0.53%  "

$ws.Range("E31").Value = "  -3.79%  "

$ws.Range("E32").Value = "  +2.60%  "

$ws.Range("D33").Value = "'0.05304"
$ws.Range("E33").Value = "  +3.32%  "

$ws.Range("D34").Value = "'1.938"
$ws.Range("E34").Value = "  +0.67%  "

$ws.Range("E35").Value = "  +1.24%  "

$ws.Range("D36").Value = "'0.7447"
$ws.Range("E36").Value = "  -7.34%  "

$ws.Range("D37").Value = "'2.696"
$ws.Range("E37").Value = "  +0.44%  "

$ws.Range("D38").Value = "'0.01869"
$ws.Range("E38").Value = "  +1.31%  "

$ws.Range("D39").Value = "1.228.61"
$ws.Range("E39").Value = "  +4.98%  "

$ws.Range("D40").Value = "'2.729"
$ws.Range("E40").Value = "  +1.10%  "

$ws.Range("D41").Value = "'6.518"
$ws.Range("E41").Value = "  +5.37%  "

$ws.Range("D42").Value = "'110.90"
$ws.Range("E42").Value = "  +8.92%  "

$ws.Range("D43").Value = "'0.8918"
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Value = "'72.95"
$ws.Range("E44").Value = "  +0.47%  "

$ws.Range("D45").Value = "'0.9999"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").Value = "2.020.46"
$ws.Range("E46").Value = "  +0.25%  "

$ws.Range("D47").Value = "'1.813"
$ws.Range("E47").Value = "  +1.64%  "

$ws.Range("D48").Value = "'0.5214"
$ws.Range("E48").Value = "  +0.67%  "

$ws.Range("E49").Value = "  +2.96%  "

$ws.Range("D50").Value = "'9.433"
$ws.Range("E50").Value = "  +2.37%  "

$ws.Range("E51").Value = "  +1.50%  "
